# Add a new "2022-Q3" sheet right before the existing "2022-Q2" sheet,
# containing the quarterly fund-holdings detail table, and insert a
# corresponding new summary row ("2022-Q3", 2, 0.07) at the top of the
# "总计" sheet's data table (shifting the older quarters down by one row).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- 1. Duplicate the "2022-Q2" sheet (this clones all formatting/styles)
#        and rename/trim it down into the new "2022-Q3" sheet. ---
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item("2022-Q2 (2)")
$q3Sheet.Name = "2022-Q3"

# Only 2 holdings this quarter -- drop the now-unneeded extra rows (4-10),
# leaving the header (row 1) + two data rows (2-3).
$q3Sheet.Range("A4:H10").Delete(-4162)

# The text-like columns (fund code / name / % figures) must stay text, not
# get auto-coerced to numbers, so force a text format before writing them,
# then strip the format back off again (the source cells carry no explicit
# style) while leaving the already-committed text value untouched.
$textCells = $q3Sheet.Range("B2:G3")
$textCells.NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "004634"
$q3Sheet.Range("C2").Value = "新疆前海联合泳涛灵活配置混合A"
$q3Sheet.Range("D2").Value = "1.20"
$q3Sheet.Range("E2").Value = "92.31"
$q3Sheet.Range("F2").Value = "4.33"
$q3Sheet.Range("G2").Value = "0.0520"
$q3Sheet.Range("H2").Value = 9

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "007041"
$q3Sheet.Range("C3").Value = "新疆前海联合泳涛灵活配置混合C"
$q3Sheet.Range("D3").Value = "0.31"
$q3Sheet.Range("E3").Value = "92.31"
$q3Sheet.Range("F3").Value = "4.33"
$q3Sheet.Range("G3").Value = "0.0134"
$q3Sheet.Range("H3").Value = 9

$textCells.ClearFormats()

# --- 2. Insert a new row at the top of the "总计" sheet's data (row 2),
#         pushing the existing quarters down by one. ---
$summary.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the row above -- clear it
# so it matches the plain (unstyled) look of the other data rows.
$summary.Range("B2:D2").ClearFormats()

# A2 needs the same "index" styling (bold/border) as the other A-column
# cells -- copy it down from A3 rather than re-building it property by
# property.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.07000000000000001

# Fix up the running index (column A) for the rows that got shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# Copying a sheet activates it; restore the originally-active sheet so the
# workbook's selection state is left exactly as it was found.
$wb.Worksheets.Item("2021-Q3").Activate()
